$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 1477.4166
$ws.Cells.Item(19, 9).Value = 607
$ws.Cells.Item(19, 11).Value = 607
$ws.Cells.Item(19, 13).Value = -432
# Row 98
$ws.Cells.Item(98, 8).Value = 3076.75
$ws.Cells.Item(98, 9).Value = 3230.5715
$ws.Cells.Item(98, 10).Value = 2000
$ws.Cells.Item(98, 11).Value = 3230.5715
$ws.Cells.Item(98, 12).Value = 2000
$ws.Cells.Item(98, 13).Value = -1732.5715
$ws.Cells.Item(98, 14).Value = -4996
# Row 122
$ws.Cells.Item(122, 8).Value = 3076.75
$ws.Cells.Item(122, 9).Value = 3230.5715
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 9691.7145
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -7241.7145
$ws.Cells.Item(122, 14).Value = -10900
# Row 137
$ws.Cells.Item(137, 8).Value = 1974.5
$ws.Cells.Item(137, 9).Value = 1638.5
$ws.Cells.Item(137, 10).Value = 2058.5
$ws.Cells.Item(137, 11).Value = 4915.5
$ws.Cells.Item(137, 12).Value = 6175.5
$ws.Cells.Item(137, 13).Value = -2365.5
$ws.Cells.Item(137, 14).Value = -11275.5
# Row 138
$ws.Cells.Item(138, 8).Value = 3452
$ws.Cells.Item(138, 10).Value = 2270.6667
$ws.Cells.Item(138, 12).Value = 6812.000100000001
$ws.Cells.Item(138, 14).Value = -17092.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5054.7144
$ws.Cells.Item(32, 9).Value = 3504.32
$ws.Cells.Item(32, 10).Value = 8930.700000000001
$ws.Cells.Item(32, 11).Value = 3504.32
$ws.Cells.Item(32, 12).Value = 8930.700000000001
$ws.Cells.Item(32, 13).Value = -3217.32
$ws.Cells.Item(32, 14).Value = -9504.700000000001
# Row 43
$ws.Cells.Item(43, 8).Value = 34350
$ws.Cells.Item(43, 10).Value = 34350
$ws.Cells.Item(43, 12).Value = 34350
$ws.Cells.Item(43, 14).Value = -34976
# Row 45
$ws.Cells.Item(45, 8).Value = 9001569
$ws.Cells.Item(45, 9).Value = 15001424
$ws.Cells.Item(45, 11).Value = 15001424
$ws.Cells.Item(45, 13).Value = -15001047
# Row 63
$ws.Cells.Item(63, 8).Value = 4027.5557
$ws.Cells.Item(63, 9).Value = 3593.5
$ws.Cells.Item(63, 10).Value = 7500
$ws.Cells.Item(63, 11).Value = 3593.5
$ws.Cells.Item(63, 12).Value = 7500
$ws.Cells.Item(63, 13).Value = -2907.5
$ws.Cells.Item(63, 14).Value = -8872
# Row 66
$ws.Cells.Item(66, 8).Value = 4027.5557
$ws.Cells.Item(66, 9).Value = 3593.5
$ws.Cells.Item(66, 10).Value = 7500
$ws.Cells.Item(66, 11).Value = 17967.5
$ws.Cells.Item(66, 12).Value = 37500
$ws.Cells.Item(66, 13).Value = -14535.5
$ws.Cells.Item(66, 14).Value = -44364
# Row 74
$ws.Cells.Item(74, 8).Value = 4000
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 14).Value = -5748
# Row 77
$ws.Cells.Item(77, 8).Value = 4000
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 14).Value = -28736
# Row 97
$ws.Cells.Item(97, 8).Value = 276.85715
$ws.Cells.Item(97, 9).Value = 210
$ws.Cells.Item(97, 11).Value = 210
$ws.Cells.Item(97, 13).Value = 286
# Row 122
$ws.Cells.Item(122, 8).Value = 1874.381
$ws.Cells.Item(122, 9).Value = 1973.5333
$ws.Cells.Item(122, 11).Value = 5920.5999
$ws.Cells.Item(122, 13).Value = -3470.5999
# Row 132
$ws.Cells.Item(132, 8).Value = 2899.5
$ws.Cells.Item(132, 9).Value = 1999
$ws.Cells.Item(132, 10).Value = 3124.625
$ws.Cells.Item(132, 11).Value = 5997
$ws.Cells.Item(132, 12).Value = 9373.875
$ws.Cells.Item(132, 13).Value = -3467
$ws.Cells.Item(132, 14).Value = -14433.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 202070
$ws.Cells.Item(86, 9).Value = 2214.2856
$ws.Cells.Item(86, 11).Value = 2214.2856
$ws.Cells.Item(86, 13).Value = -1091.2856
# Row 89
$ws.Cells.Item(89, 8).Value = 202070
$ws.Cells.Item(89, 9).Value = 2214.2856
$ws.Cells.Item(89, 11).Value = 11071.428
$ws.Cells.Item(89, 13).Value = -5455.428
# Row 94
$ws.Cells.Item(94, 8).Value = 656.375
$ws.Cells.Item(94, 9).Value = 482.75
$ws.Cells.Item(94, 11).Value = 482.75
$ws.Cells.Item(94, 13).Value = -31.75
# Row 99
$ws.Cells.Item(99, 8).Value = 1377.2667
$ws.Cells.Item(99, 9).Value = 1286.2
$ws.Cells.Item(99, 11).Value = 1286.2
$ws.Cells.Item(99, 13).Value = 211.8
# Row 105
$ws.Cells.Item(105, 8).Value = 2111.923
$ws.Cells.Item(105, 9).Value = 2148.6191
$ws.Cells.Item(105, 11).Value = 2148.6191
$ws.Cells.Item(105, 13).Value = -401.6190999999999
# Row 134
$ws.Cells.Item(134, 8).Value = 13804.934
$ws.Cells.Item(134, 9).Value = 15769.083
$ws.Cells.Item(134, 11).Value = 47307.249
$ws.Cells.Item(134, 13).Value = -44772.249

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5140
$ws.Cells.Item(31, 9).Value = 2500
$ws.Cells.Item(31, 10).Value = 6900
$ws.Cells.Item(31, 11).Value = 2500
$ws.Cells.Item(31, 12).Value = 6900
$ws.Cells.Item(31, 13).Value = -2205
$ws.Cells.Item(31, 14).Value = -7490
# Row 34
$ws.Cells.Item(34, 8).Value = 5140
$ws.Cells.Item(34, 9).Value = 2500
$ws.Cells.Item(34, 10).Value = 6900
$ws.Cells.Item(34, 11).Value = 2500
$ws.Cells.Item(34, 12).Value = 6900
$ws.Cells.Item(34, 13).Value = -2298
$ws.Cells.Item(34, 14).Value = -7304
# Row 58
$ws.Cells.Item(58, 8).Value = 1299.6666
$ws.Cells.Item(58, 9).Value = 1449.5
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 1449.5
$ws.Cells.Item(58, 12).Value = 1000
$ws.Cells.Item(58, 13).Value = -1246.5
$ws.Cells.Item(58, 14).Value = -1406
# Row 122
$ws.Cells.Item(122, 8).Value = 1873.65
$ws.Cells.Item(122, 9).Value = 1721
$ws.Cells.Item(122, 10).Value = 2060.2222
$ws.Cells.Item(122, 11).Value = 5163
$ws.Cells.Item(122, 12).Value = 6180.6666
$ws.Cells.Item(122, 13).Value = -2713
$ws.Cells.Item(122, 14).Value = -11080.6666
# Row 136
$ws.Cells.Item(136, 8).Value = 1299.6666
$ws.Cells.Item(136, 9).Value = 1449.5
$ws.Cells.Item(136, 10).Value = 1000
$ws.Cells.Item(136, 11).Value = 4348.5
$ws.Cells.Item(136, 12).Value = 3000
$ws.Cells.Item(136, 13).Value = -1798.5
$ws.Cells.Item(136, 14).Value = -8100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Cells.Item(98, 8).Value = 538.2
# Row 122
$ws.Cells.Item(122, 8).Value = 1494.7142
$ws.Cells.Item(122, 9).Value = 1092.3334
$ws.Cells.Item(122, 11).Value = 9831.000599999999
$ws.Cells.Item(122, 13).Value = -7381.000599999999
# Row 125
$ws.Cells.Item(125, 8).Value = 3780
$ws.Cells.Item(125, 9).Value = 1652
$ws.Cells.Item(125, 11).Value = 4956
$ws.Cells.Item(125, 13).Value = -36
# Row 132
$ws.Cells.Item(132, 8).Value = 1643.8889
$ws.Cells.Item(132, 10).Value = 2097.5
$ws.Cells.Item(132, 12).Value = 18877.5
$ws.Cells.Item(132, 14).Value = -23937.5
# Row 139
$ws.Cells.Item(139, 8).Value = 7987.533
$ws.Cells.Item(139, 9).Value = 8422.357
$ws.Cells.Item(139, 11).Value = 25267.071
$ws.Cells.Item(139, 13).Value = -20127.071

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 4266.909
$ws.Cells.Item(70, 9).Value = 4654.5
$ws.Cells.Item(70, 10).Value = 3233.3333
$ws.Cells.Item(70, 11).Value = 4654.5
$ws.Cells.Item(70, 12).Value = 3233.3333
$ws.Cells.Item(70, 13).Value = -4384.5
$ws.Cells.Item(70, 14).Value = -3773.3333
# Row 73
$ws.Cells.Item(73, 8).Value = 4266.909
$ws.Cells.Item(73, 9).Value = 4654.5
$ws.Cells.Item(73, 10).Value = 3233.3333
$ws.Cells.Item(73, 11).Value = 4654.5
$ws.Cells.Item(73, 12).Value = 3233.3333
$ws.Cells.Item(73, 13).Value = -3718.5
$ws.Cells.Item(73, 14).Value = -5105.3333
# Row 132
$ws.Cells.Item(132, 8).Value = 3783.625
$ws.Cells.Item(132, 9).Value = 3401.25
$ws.Cells.Item(132, 10).Value = 4166
$ws.Cells.Item(132, 11).Value = 10203.75
$ws.Cells.Item(132, 12).Value = 12498
$ws.Cells.Item(132, 13).Value = -7673.75
$ws.Cells.Item(132, 14).Value = -17558
# Row 134
$ws.Cells.Item(134, 8).Value = 44285.285
$ws.Cells.Item(134, 10).Value = 44285.285
$ws.Cells.Item(134, 12).Value = 132855.855
$ws.Cells.Item(134, 14).Value = -137925.855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 6668.8125
$ws.Cells.Item(40, 9).Value = 3275.25
$ws.Cells.Item(40, 11).Value = 3275.25
$ws.Cells.Item(40, 13).Value = -3139.25
# Row 46
$ws.Cells.Item(46, 8).Value = 1501.8462
$ws.Cells.Item(46, 10).Value = 1608.25
$ws.Cells.Item(46, 12).Value = 1608.25
$ws.Cells.Item(46, 14).Value = -1984.25
# Row 93
$ws.Cells.Item(93, 8).Value = 835.5238000000001
$ws.Cells.Item(93, 9).Value = 742.3333
$ws.Cells.Item(93, 10).Value = 1394.6666
$ws.Cells.Item(93, 11).Value = 742.3333
$ws.Cells.Item(93, 12).Value = 1394.6666
$ws.Cells.Item(93, 13).Value = 505.6667
$ws.Cells.Item(93, 14).Value = -3890.6666
# Row 122
$ws.Cells.Item(122, 8).Value = 9454.546
$ws.Cells.Item(122, 9).Value = 10000
$ws.Cells.Item(122, 11).Value = 30000
$ws.Cells.Item(122, 13).Value = -27550
# Row 132
$ws.Cells.Item(132, 8).Value = 1933
$ws.Cells.Item(132, 9).Value = 1319.7
$ws.Cells.Item(132, 10).Value = 4999.5
$ws.Cells.Item(132, 11).Value = 3959.1
$ws.Cells.Item(132, 12).Value = 14998.5
$ws.Cells.Item(132, 13).Value = -1429.1
$ws.Cells.Item(132, 14).Value = -20058.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 829.2857
$ws.Cells.Item(107, 9).Value = 693.7273
$ws.Cells.Item(107, 11).Value = 2081.1819
$ws.Cells.Item(107, 13).Value = -161.1819
# Row 109
$ws.Cells.Item(109, 8).Value = 42358.332
$ws.Cells.Item(109, 10).Value = 42358.332
$ws.Cells.Item(109, 12).Value = 42358.332
$ws.Cells.Item(109, 14).Value = -45132.332
# Row 122
$ws.Cells.Item(122, 8).Value = 260592.33
$ws.Cells.Item(122, 9).Value = 260592.33
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 781776.99
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -779326.99
$ws.Cells.Item(122, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value = 4917.8335
$ws.Cells.Item(126, 9).Value = 4376.5625
$ws.Cells.Item(126, 10).Value = 6000.375
$ws.Cells.Item(126, 11).Value = 13129.6875
$ws.Cells.Item(126, 12).Value = 18001.125
$ws.Cells.Item(126, 13).Value = -10659.6875
$ws.Cells.Item(126, 14).Value = -22941.125
# Row 132
$ws.Cells.Item(132, 8).Value = 1710.95
$ws.Cells.Item(132, 9).Value = 954.73334
$ws.Cells.Item(132, 11).Value = 2864.20002
$ws.Cells.Item(132, 13).Value = -334.2000200000002
